$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create the hidden ValidationData sheet, placed after the main sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ValidationData"

# --- Column data for the validation lists ---
$Icol = @("1K-10K", "10K-25K", "25K-50K", "50K-100K", "100K-250K", "250K-500K", "500K-1M", "1M+")
$Jcol = @("Pre-seed", "Seed", "Series A", "Series B", "Series C", "Growth", "All")
$Kcol = @("-- Tech`n    B2B SaaS", "Fintech", "Healthtech", "AI/ML", "Deep tech", "Climate tech", "Consumer", "E-commerce", "Marketplace", "Gaming", "Web3", "Developer tools", "Cybersecurity", "Logistics", "AdTech", "PropTech", "InsurTech", "-- Non-Tech / Other`n    Agriculture", "Automotive", "Biotechnology", "Construction", "Consulting", "Consumer Goods", "Education", "Energy", "Entertainment", "Environmental Services", "Fashion", "Food & Beverage", "Government", "Healthcare Services", "Hospitality", "Human Resources", "Insurance", "Legal", "Manufacturing", "Media", "Non-profit", "Pharmaceuticals", "Real Estate", "Retail", "Telecommunications", "Transportation", "Utilities", "Other")
$Lcol = @("Global", "North America", "South America", "LATAM", "Europe", "Western Europe", "Eastern Europe", "Continental Europe", "Middle East", "Africa", "Asia", "East Asia", "South Asia", "South East Asia", "Oceania", "EMEA", "Emerging Markets")
$Mcol = @("hands-on", "passive", "advisory", "network-focused")
$Pcol = @("1-3 days", "1 week", "2 weeks", "1 month", "2+ months")
$Qcol = @("form", "email", "other")
$Tcol = @("simple", "standard", "comprehensive")
$Ucol = @("pitch_deck", "video", "financial_projections", "business_plan", "traction_data")
$Zcol = @("FREE", "PRO", "MAX")

# --- Write cells row-by-row, columns in I,J,K,L,M,P,Q,T,U,Z order, so the
# shared-string table is interned in the same order as the target workbook ---
$ws2.Range("I1").Value = $Icol[0]
$ws2.Range("J1").Value = $Jcol[0]
$ws2.Range("K1").Value = $Kcol[0]
$ws2.Range("L1").Value = $Lcol[0]
$ws2.Range("M1").Value = $Mcol[0]
$ws2.Range("P1").Value = $Pcol[0]
$ws2.Range("Q1").Value = $Qcol[0]
$ws2.Range("T1").Value = $Tcol[0]
$ws2.Range("U1").Value = $Ucol[0]
$ws2.Range("Z1").Value = $Zcol[0]
$ws2.Range("I2").Value = $Icol[1]
$ws2.Range("J2").Value = $Jcol[1]
$ws2.Range("K2").Value = $Kcol[1]
$ws2.Range("L2").Value = $Lcol[1]
$ws2.Range("M2").Value = $Mcol[1]
$ws2.Range("P2").Value = $Pcol[1]
$ws2.Range("Q2").Value = $Qcol[1]
$ws2.Range("T2").Value = $Tcol[1]
$ws2.Range("U2").Value = $Ucol[1]
$ws2.Range("Z2").Value = $Zcol[1]
$ws2.Range("I3").Value = $Icol[2]
$ws2.Range("J3").Value = $Jcol[2]
$ws2.Range("K3").Value = $Kcol[2]
$ws2.Range("L3").Value = $Lcol[2]
$ws2.Range("M3").Value = $Mcol[2]
$ws2.Range("P3").Value = $Pcol[2]
$ws2.Range("Q3").Value = $Qcol[2]
$ws2.Range("T3").Value = $Tcol[2]
$ws2.Range("U3").Value = $Ucol[2]
$ws2.Range("Z3").Value = $Zcol[2]
$ws2.Range("I4").Value = $Icol[3]
$ws2.Range("J4").Value = $Jcol[3]
$ws2.Range("K4").Value = $Kcol[3]
$ws2.Range("L4").Value = $Lcol[3]
$ws2.Range("M4").Value = $Mcol[3]
$ws2.Range("P4").Value = $Pcol[3]
$ws2.Range("U4").Value = $Ucol[3]
$ws2.Range("I5").Value = $Icol[4]
$ws2.Range("J5").Value = $Jcol[4]
$ws2.Range("K5").Value = $Kcol[4]
$ws2.Range("L5").Value = $Lcol[4]
$ws2.Range("P5").Value = $Pcol[4]
$ws2.Range("U5").Value = $Ucol[4]
$ws2.Range("I6").Value = $Icol[5]
$ws2.Range("J6").Value = $Jcol[5]
$ws2.Range("K6").Value = $Kcol[5]
$ws2.Range("L6").Value = $Lcol[5]
$ws2.Range("I7").Value = $Icol[6]
$ws2.Range("J7").Value = $Jcol[6]
$ws2.Range("K7").Value = $Kcol[6]
$ws2.Range("L7").Value = $Lcol[6]
$ws2.Range("I8").Value = $Icol[7]
$ws2.Range("K8").Value = $Kcol[7]
$ws2.Range("L8").Value = $Lcol[7]
$ws2.Range("K9").Value = $Kcol[8]
$ws2.Range("L9").Value = $Lcol[8]
$ws2.Range("K10").Value = $Kcol[9]
$ws2.Range("L10").Value = $Lcol[9]
$ws2.Range("K11").Value = $Kcol[10]
$ws2.Range("L11").Value = $Lcol[10]
$ws2.Range("K12").Value = $Kcol[11]
$ws2.Range("L12").Value = $Lcol[11]
$ws2.Range("K13").Value = $Kcol[12]
$ws2.Range("L13").Value = $Lcol[12]
$ws2.Range("K14").Value = $Kcol[13]
$ws2.Range("L14").Value = $Lcol[13]
$ws2.Range("K15").Value = $Kcol[14]
$ws2.Range("L15").Value = $Lcol[14]
$ws2.Range("K16").Value = $Kcol[15]
$ws2.Range("L16").Value = $Lcol[15]
$ws2.Range("K17").Value = $Kcol[16]
$ws2.Range("L17").Value = $Lcol[16]
$ws2.Range("K18").Value = $Kcol[17]
$ws2.Range("K19").Value = $Kcol[18]
$ws2.Range("K20").Value = $Kcol[19]
$ws2.Range("K21").Value = $Kcol[20]
$ws2.Range("K22").Value = $Kcol[21]
$ws2.Range("K23").Value = $Kcol[22]
$ws2.Range("K24").Value = $Kcol[23]
$ws2.Range("K25").Value = $Kcol[24]
$ws2.Range("K26").Value = $Kcol[25]
$ws2.Range("K27").Value = $Kcol[26]
$ws2.Range("K28").Value = $Kcol[27]
$ws2.Range("K29").Value = $Kcol[28]
$ws2.Range("K30").Value = $Kcol[29]
$ws2.Range("K31").Value = $Kcol[30]
$ws2.Range("K32").Value = $Kcol[31]
$ws2.Range("K33").Value = $Kcol[32]
$ws2.Range("K34").Value = $Kcol[33]
$ws2.Range("K35").Value = $Kcol[34]
$ws2.Range("K36").Value = $Kcol[35]
$ws2.Range("K37").Value = $Kcol[36]
$ws2.Range("K38").Value = $Kcol[37]
$ws2.Range("K39").Value = $Kcol[38]
$ws2.Range("K40").Value = $Kcol[39]
$ws2.Range("K41").Value = $Kcol[40]
$ws2.Range("K42").Value = $Kcol[41]
$ws2.Range("K43").Value = $Kcol[42]
$ws2.Range("K44").Value = $Kcol[43]
$ws2.Range("K45").Value = $Kcol[44]

# --- Update data validations on sheet1 to reference ValidationData ranges ---
$v = $ws1.Range("I10:I1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$I`$1:`$I`$8")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("I2:I1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$I`$1:`$I`$8")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("J10:J1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$J`$1:`$J`$7")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("J2:J1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$J`$1:`$J`$7")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("K10:K1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$K`$1:`$K`$45")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("K2:K1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$K`$1:`$K`$45")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("L10:L1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$L`$1:`$L`$17")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("L2:L1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$L`$1:`$L`$17")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("M10:M1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$M`$1:`$M`$4")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("M2:M1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$M`$1:`$M`$4")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("P10:P1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$P`$1:`$P`$5")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("P2:P1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$P`$1:`$P`$5")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("Q10:Q1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$Q`$1:`$Q`$3")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("Q2:Q1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$Q`$1:`$Q`$3")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("T10:T1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$T`$1:`$T`$3")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("T2:T1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$T`$1:`$T`$3")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("U10:U1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$U`$1:`$U`$5")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("U2:U1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$U`$1:`$U`$5")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("Z10:Z1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$Z`$1:`$Z`$3")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."
$v = $ws1.Range("Z2:Z1000").Validation
$v.Modify(3, 1, 1, "ValidationData!`$Z`$1:`$Z`$3")
$v.ErrorTitle = "Invalid Value"
$v.ErrorMessage = "Please select a value from the dropdown list."

# --- Bold the header row on sheet1 ---
$ws1.Range("A1:Z1").Font.Bold = $true

# --- Hide the ValidationData sheet ---
$ws2.Visible = $false

# --- Keep the main sheet as the active/selected one ---
$ws1.Activate()

